$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.777.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.29%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.313.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.11%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'539.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.15%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'132.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.75%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.06%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  +2.05%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'2.309.12"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.20%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  -1.52%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  -0.95%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  +0.81%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  -0.60%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'23.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.14%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'2.722.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.00%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'58.686.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.38%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  -0.46%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.277.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.73%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'10.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.23%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'4.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.05%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'316.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.13%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'6.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.81%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  +0.06%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'63.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.34%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.173"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.26%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.09%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'7.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.57%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'1.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.77%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("B29").Value = "'Monero"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'171.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.57%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("B30").Value = "'PancakeSwap"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'1.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.99%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.0₃0735"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.05%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'1.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.64%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'5.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.61%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  +0.67%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'17.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.62%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  -0.02%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  +0.05%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'1.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.68%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'4.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.81%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  +0.08%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'298.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.86%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'141.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.01%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  +0.10%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.0957"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.01%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = "'  -1.54%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  -0.70%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'18.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.87%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  -2.71%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'10.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.23%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E51").Value = "'  +0.63%  "
$ws.Range("E51").Style = "Normal"
